$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = 8.802599999999998
    "B6"  = 6.417100000000002
    "B7"  = 4.964900000000003
    "B8"  = 6.408700000000002
    "B16" = 6.682199999999999
    "B20" = 9.435699999999995
    "B21" = 9.436999999999994
    "B28" = 5.981299999999999
    "B29" = 5.185300000000004
    "B30" = 5.495500000000003
    "B32" = 7.33649999999999
    "B40" = 8.985099999999994
    "B46" = 6.112200000000001
    "B51" = 5.572900000000001
    "B52" = 5.266099999999999
    "B57" = 5.359699999999997
    "B59" = 4.873299999999998
    "B62" = 6.0513
    "B66" = 5.405600000000001
    "B73" = 8.410099999999998
    "B74" = 9.153599999999996
    "B77" = 8.891700000000005
    "B92" = 4.845199999999998
    "B100" = 5.429699999999996
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
